$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Power Query refresh results: updated "월별 누적별풍선" totals (col C) and
# "새로고침시간" refresh timestamp (col D) for every data row (rows 2-12).

$ws.Range("C2").Value = 747963
$ws.Range("D2").Value = 46016.947262083333

$ws.Range("D3").Value = 46016.947262083333

$ws.Range("C4").Value = 534610
$ws.Range("D4").Value = 46016.947262083333

$ws.Range("C5").Value = 527910
$ws.Range("D5").Value = 46016.947262083333

$ws.Range("C6").Value = 428446
$ws.Range("D6").Value = 46016.947262083333

$ws.Range("C7").Value = 409390
$ws.Range("D7").Value = 46016.947262083333

$ws.Range("D8").Value = 46016.947262083333

$ws.Range("C9").Value = 286854
$ws.Range("D9").Value = 46016.947262083333

$ws.Range("C10").Value = 202994
$ws.Range("D10").Value = 46016.947262083333

$ws.Range("C11").Value = 169099
$ws.Range("D11").Value = 46016.947262083333

$ws.Range("D12").Value = 46016.947262083333

# Update the active selection left behind when the file was last saved.
$ws.Range("E14").Select() | Out-Null
